# Fruta / hortaliza, semanal
# The data rows (2-37) get reshuffled: each row's full record (date, variety,
# quality, volume, prices, unit, origin, price/kg, etc.) moves to a different
# row position. Column count: A..R (1..18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 37
$firstCol = 1
$lastCol = 18

# Mapping: source (current/"before") row -> destination (target/"after") row.
$rowMap = @{
    2 = 35
    3 = 7
    4 = 3
    5 = 14
    6 = 17
    7 = 24
    8 = 27
    9 = 8
    10 = 23
    11 = 31
    12 = 26
    13 = 30
    14 = 11
    15 = 15
    16 = 33
    17 = 9
    18 = 32
    19 = 20
    20 = 12
    21 = 18
    22 = 29
    23 = 28
    24 = 10
    25 = 13
    26 = 25
    27 = 6
    28 = 2
    29 = 36
    30 = 34
    31 = 19
    32 = 5
    33 = 22
    34 = 37
    35 = 16
    36 = 4
    37 = 21
}

# 1) Snapshot every source row's values before any writes happen, so that
#    writing to destination rows never clobbers data we still need to read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each snapshotted row out to its mapped destination row.
foreach ($srcRow in $rowMap.Keys) {
    $dstRow = $rowMap[$srcRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($dstRow, $c).Value2 = $rowVals[$c - 1]
    }
}
